$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update final score (column K) values
$ws.Range("K2").Value = 60.9
$ws.Range("K3").Value = 57.1
$ws.Range("K4").Value = 55.7
$ws.Range("K5").Value = 51.1
$ws.Range("K6").Value = 50.9
$ws.Range("K7").Value = 44.7

# Update MACRO_SCORE (column N) values
$ws.Range("N2").Value = 85.82376350509293
$ws.Range("N3").Value = 85.82376350509293
$ws.Range("N4").Value = 85.82376350509293
$ws.Range("N5").Value = 85.82376350509293
$ws.Range("N6").Value = 85.82376350509293
$ws.Range("N7").Value = 85.82376350509293
